$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Resource Utilization" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = 11.5
$ws1.Range("B3").Value = 0.9399999999999999

# --- Sheet 2: "Activity Times" ---
$ws2 = $wb.Worksheets.Item(2)

# Row 2 ("5.5.13 Real Property-Monthly Reviews-org" / Process)
$ws2.Range("C2").Value = 1080
$ws2.Range("D2").Value = 3
$ws2.Range("F2").Value = 810
$ws2.Range("G2").Value = 122.38

# Row 3 ("Monthly" / Start)
$ws2.Range("C3").Value = 1080

# Row 4 ("Review AM using Asset Change Tracker (5.5.13.1)" / Activity Step)
$ws2.Range("C4").Value = 1080
$ws2.Range("D4").Value = 66
$ws2.Range("F4").Value = 810
$ws2.Range("G4").Value = 35.02

# Row 5 ("Complete /Accurate?" / Gateway)
$ws2.Range("C5").Value = 66
$ws2.Range("D5").Value = 66

# Row 6 — becomes "Work with REO RPO to Correct (5.5.13.3)" / Activity Step
$ws2.Range("A6").Value = "Work with REO RPO to Correct (5.5.13.3)"
$ws2.Range("B6").Value = "Activity Step"
$ws2.Range("C6").Value = 17
$ws2.Range("D6").Value = 3
$ws2.Range("E6").Value = 86
$ws2.Range("F6").Value = 166
$ws2.Range("G6").Value = 128.67
$ws2.Range("H6:K6").ClearContents()

# Row 7 — becomes "Note Accuracy in Asset Change Tracker (5.5.13.2)" / Activity Step
$ws2.Range("A7").Value = "Note Accuracy in Asset Change Tracker (5.5.13.2)"
$ws2.Range("B7").Value = "Activity Step"
$ws2.Range("C7").Value = 3
$ws2.Range("D7").Value = 3
$ws2.Range("E7").Value = 90
$ws2.Range("F7").Value = 168
$ws2.Range("G7").Value = 132.33
$ws2.Range("H7:K7").ClearContents()

# Row 8 — becomes "Create/Post Journal Entries (5.5.13.4)" / Stop
$ws2.Range("A8").Value = "Create/Post Journal Entries (5.5.13.4)"
$ws2.Range("B8").Value = "Stop"
$ws2.Range("C8").Value = 3
$ws2.Range("D8").Value = 3
$ws2.Range("E8").Value = 95
$ws2.Range("F8").Value = 174
$ws2.Range("G8").Value = 137.67
$ws2.Range("H8:K8").ClearContents()
